# This script re-derives the T+4/T+5 (5_naive_qoq_error_series_matched_to_ifoCAST)
# "since 2021" matched-error table for a full ifoCAST evaluation: every quarter
# row drops its oldest (Q0) error and gains one additional quarter of realised
# error where available, i.e. each data row shifts one column to the left and the
# rows that have run out of realised quarters simply get shorter (their trailing
# cells are cleared instead of refilled).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B..K (indices 0..9), one array per data row (rows 2..20).
# $null marks a cell that must be cleared because the row is now shorter than before.
$newData = @(
    @(2.03161760340632, -1.146501396593679, -1.441830396593679, 1.072614603406321, 0.1428926034063205, 0.3779006034063205, 0.07076460340632051, 0.5814406034063205, 0.2940736034063205, 0.3220726034063205),
    @(-3.294071868393072, -3.589400868393072, -1.074955868393072, -2.004677868393072, -1.769669868393072, -2.076805868393072, -1.566129868393072, -1.853496868393072, -1.825497868393072, -1.732849868393072),
    @(-0.1493708818697419, 2.365074118130258, 1.435352118130258, 1.670360118130258, 1.363224118130258, 1.873900118130258, 1.586533118130258, 1.614532118130258, 1.707180118130258, 1.920156118130258),
    @(2.426248296544459, 1.496526296544459, 1.731534296544459, 1.424398296544459, 1.935074296544459, 1.647707296544459, 1.675706296544459, 1.768354296544459, 1.981330296544459, 1.480670296544459),
    @(-0.5449296406117954, -0.3099216406117954, -0.6170576406117954, -0.1063816406117954, -0.3937486406117954, -0.3657496406117954, -0.2731016406117954, -0.06012564061179543, -0.5607856406117955, -0.2710286406117954),
    @(0.2726728964218456, -0.03446310357815441, 0.4762128964218456, 0.1888458964218456, 0.2168448964218456, 0.3094928964218456, 0.5224688964218456, 0.02180889642184558, 0.3115658964218456, 0.05367389642184558),
    @(-0.1301381443875122, 0.3805378556124878, 0.0931708556124878, 0.1211698556124878, 0.2138178556124878, 0.4267938556124878, -0.07386614438751221, 0.2158908556124878, -0.04200114438751221, 0.3403798556124878),
    @(0.6409568926112106, 0.3535898926112106, 0.3815888926112106, 0.4742368926112106, 0.6872128926112107, 0.1865528926112106, 0.4763098926112106, 0.2184178926112106, 0.6007988926112107, 0.1895217986112106),
    @(-0.4818651035472806, -0.4538661035472806, -0.3612181035472806, -0.1482421035472806, -0.6489021035472806, -0.3591451035472806, -0.6170371035472806, -0.2346561035472806, -0.6459331975472806, -0.3589771035472806),
    @(-0.3537867436446591, -0.2611387436446591, -0.04816274364465911, -0.5488227436446591, -0.2590657436446591, -0.5169577436446591, -0.1345767436446591, -0.5458538376446591, -0.2588977436446591, $null),
    @(0.09264776243503714, 0.3056237624350371, -0.1950362375649629, 0.09472076243503715, -0.1631712375649629, 0.2192097624350371, -0.1920673315649629, 0.09488876243503713, $null, $null),
    @(0.1574463720025918, -0.3432136279974082, -0.05345662799740819, -0.3113486279974082, 0.0710323720025918, -0.3402447219974082, -0.0532886279974082, $null, $null, $null),
    @(-0.5006596170015631, -0.2109026170015632, -0.4687946170015632, -0.08641361700156319, -0.4976907110015631, -0.2107346170015632, $null, $null, $null, $null),
    @(0.2897568395245076, 0.03186483952450761, 0.4142458395245076, 0.002968745524507627, 0.2899248395245076, $null, $null, $null, $null, $null),
    @(-0.2484930410109615, 0.1338879589890384, -0.2773891350109615, 0.009566958989038449, $null, $null, $null, $null, $null, $null),
    @(0.3058628168340501, -0.1054142771659499, 0.1815418168340501, $null, $null, $null, $null, $null, $null, $null),
    @(-0.4112768169122814, -0.1243207229122814, $null, $null, $null, $null, $null, $null, $null, $null),
    @(0.4880092297750048, $null, $null, $null, $null, $null, $null, $null, $null, $null),
    @($null, $null, $null, $null, $null, $null, $null, $null, $null, $null)
)

for ($i = 0; $i -lt $newData.Length; $i++) {
    $rowNum = $i + 2         # data starts on row 2
    $rowVals = $newData[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $colNum = $j + 2     # column B = 2
        $cell = $ws.Cells.Item($rowNum, $colNum)
        if ($null -eq $rowVals[$j]) {
            $cell.ClearContents()
        } else {
            $cell.Value = $rowVals[$j]
        }
    }
}
